$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "69.781.26"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.40%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "3.380.27"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.23%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "582.69"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.38%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "180.24"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.82%  "
$ws.Cells.Item(7, 5).Value = "  +0.09%  "
$ws.Cells.Item(8, 5).Value = "  +0.48%  "
$ws.Cells.Item(9, 5).Value = "  +8.97%  "
$ws.Cells.Item(10, 5).Value = "  +1.23%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "48.55"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.14%  "
$ws.Cells.Item(12, 5).Value = "  +4.39%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "683.46"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.64%  "
$ws.Cells.Item(14, 5).Value = "  +2.14%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "3.926.15"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.13%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "69.752.51"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.31%  "
$ws.Cells.Item(17, 5).Value = "  +0.92%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "3.386.42"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.26%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "17.71"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.34%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "11.30"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.11%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.913"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +1.94%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "17.36"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +2.44%  "
$ws.Cells.Item(23, 5).Value = "  -1.56%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "101.87"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.30%  "
$ws.Cells.Item(25, 5).Value = "  -0.30%  "
$ws.Cells.Item(26, 5).Value = "  +0.31%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "9.78"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +3.05%  "
$ws.Cells.Item(28, 5).Value = "  +1.55%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "8.77"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.64%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "3.84"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +15.33%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "11.09"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.33%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "556.09"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.84%  "
$ws.Cells.Item(34, 5).Value = "  +0.75%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "57.99"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.79%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.08%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "3.608.26"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.87%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.139"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.46%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "35.42"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.23%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0732"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +8.61%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "2.77"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +5.58%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "3.33"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +4.87%  "
$ws.Cells.Item(43, 2).Value = "VeChain"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.0428"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +3.77%  "
$ws.Cells.Item(44, 2).Value = "TheGraph"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.338"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.83%  "
$ws.Cells.Item(45, 2).Value = "ThetaToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "2.68"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.93%  "
$ws.Cells.Item(46, 2).Value = "Stellar"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.129"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.08%  "
$ws.Cells.Item(47, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -0.12%  "
$ws.Cells.Item(48, 2).Value = "Mantle"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.38"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +3.71%  "
$ws.Cells.Item(49, 2).Value = "Monero"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "130.62"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -0.05%  "
$ws.Cells.Item(50, 2).Value = "CoreDAO"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "2.62"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.18%  "
$ws.Cells.Item(51, 2).Value = "THORChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "7.50"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.82%  "
